$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 42, shifting existing data (old rows 42-84) down to 44-86
$ws.Rows.Item(42).Insert()
$ws.Rows.Item(42).Insert()

# Populate new row 42
$ws.Cells.Item(42,1).Value = 2
$ws.Cells.Item(42,2).Value = 'Comercializadora del Agro de Limarí'
$ws.Cells.Item(42,3).Value = 'Coquimbo'
$ws.Cells.Item(42,4).Value = 44601
$ws.Cells.Item(42,5).Value = 4
$ws.Cells.Item(42,6).Value = 'Fruta'
$ws.Cells.Item(42,7).Value = 100103
$ws.Cells.Item(42,8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(42,9).Value = 100103004
$ws.Cells.Item(42,10).Value = 'Durazno'
$ws.Cells.Item(42,11).Value = 'Loadel'
$ws.Cells.Item(42,12).Value = 'Especial'
$ws.Cells.Item(42,13).Value = 100
$ws.Cells.Item(42,14).Value = 22500
$ws.Cells.Item(42,15).Value = 23000
$ws.Cells.Item(42,16).Value = 22750
$ws.Cells.Item(42,17).Value = '$/caja 16 kilos empedrada'
$ws.Cells.Item(42,18).Value = 'Región de O''Higgins'
$ws.Cells.Item(42,19).Value = 1422
$ws.Cells.Item(42,20).Value = 16

# Populate new row 43
$ws.Cells.Item(43,1).Value = 2
$ws.Cells.Item(43,2).Value = 'Comercializadora del Agro de Limarí'
$ws.Cells.Item(43,3).Value = 'Coquimbo'
$ws.Cells.Item(43,4).Value = 44601
$ws.Cells.Item(43,5).Value = 4
$ws.Cells.Item(43,6).Value = 'Fruta'
$ws.Cells.Item(43,7).Value = 100103
$ws.Cells.Item(43,8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(43,9).Value = 100103004
$ws.Cells.Item(43,10).Value = 'Durazno'
$ws.Cells.Item(43,11).Value = 'Loadel'
$ws.Cells.Item(43,12).Value = 'Primera'
$ws.Cells.Item(43,13).Value = 40
$ws.Cells.Item(43,14).Value = 18500
$ws.Cells.Item(43,15).Value = 19000
$ws.Cells.Item(43,16).Value = 18750
$ws.Cells.Item(43,17).Value = '$/caja 16 kilos empedrada'
$ws.Cells.Item(43,18).Value = 'Región de O''Higgins'
$ws.Cells.Item(43,19).Value = 1172
$ws.Cells.Item(43,20).Value = 16
